# Block API for Line Items
# Adds a new "Staging" environment column (D) to the credentials sheet,
# bumps the Quantity value in B8, and updates the customer-name block
# (row 30) to add a "VP" short-code column and capitalize the tata -> Tata
# sample value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add the new "Staging" hyperlink column (D1) ---------------
$ws.Range("D1").Value = "https://app.calcuquote.com/Staging/"
$ws.Hyperlinks.Add($ws.Range("D1"), "https://app.calcuquote.com/Staging/")
$ws.Range("D1").Style = $ws.Range("C1").Style

# --- Row 30: add VP short-code column -----------------------------------
$ws.Range("D30").Value = "VP"

# --- Row 8: bump the quantity value from 1 -> 2 ------------------------
$ws.Range("B8").Value = "2"

# --- Row 30: capitalize sample customer name ("tata" -> "Tata") -------
$ws.Range("B30").Value = "Tata"

# --- Update selection / active cell to match the edited area ----------
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B30").Select()
